# Update F-column numeric values ("想去人数" / want-to-go counts) on sheets
# "展览", "本地生活" and "全部类型" per the commit's regenerated output data.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F6").Value = 1107
$ws1.Range("F7").Value = 887
$ws1.Range("F8").Value = 271
$ws1.Range("F14").Value = 515
$ws1.Range("F18").Value = 1234
$ws1.Range("F19").Value = 2914
$ws1.Range("F20").Value = 1504
$ws1.Range("F21").Value = 738
$ws1.Range("F27").Value = 3231
$ws1.Range("F30").Value = 1441

$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F2").Value = 763

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 763
$ws4.Range("F10").Value = 1107
$ws4.Range("F11").Value = 887
$ws4.Range("F12").Value = 271
$ws4.Range("F24").Value = 515
$ws4.Range("F28").Value = 1234
$ws4.Range("F29").Value = 2914
$ws4.Range("F30").Value = 1504
$ws4.Range("F31").Value = 738
$ws4.Range("F39").Value = 3231
$ws4.Range("F42").Value = 1441
